# This workbook stores a rolling weekly price series for "Ajo" (garlic) at
# Terminal Hortofrutícola Agro Chillán in data rows 340..441 (row 1 is the
# header). A new weekly observation was added at the top of the series
# (row 340): every existing observation shifts down by one row, and the
# observation that used to be the last one (row 441) is preserved by being
# copied down into a brand-new row 442.
#
# Columns A, B, C, E, F, G and R (market id/name, region, category id,
# category name and classification) are constant for every row in this
# block, so they do not need to move. Columns D, H, I, J, K, L, M, N, O, P
# and Q hold the actual per-observation data (date, variety, quality,
# volume, prices, unit, origin, unit price, kg) and are the ones that
# shift down.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Capture the current (pre-edit) data block, rows 340-441, all columns.
$srcRange = $ws.Range("A340:R441")
$data = $srcRange.Value2

# 2) Shift the whole block down by one row: old row N (340..441) becomes
#    new row N+1 (341..442). This automatically places the former row 441
#    into the brand-new row 442.
$dstRange = $ws.Range("A341:R442")
$dstRange.Value2 = $data

# 3) Row 340 becomes the new weekly observation. The "static" columns
#    (A, B, C, E, F, G, R) keep the values already in place; only the
#    observation-specific columns change.
$ws.Cells.Item(340, 4).Value2 = 45093       # D340 - Fecha
$ws.Cells.Item(340, 8).Value2 = "Chino"     # H340 - Variedad
$ws.Cells.Item(340, 9).Value2 = "Primera"   # I340 - Calidad
$ws.Cells.Item(340, 10).Value2 = 60         # J340 - Volumen
$ws.Cells.Item(340, 11).Value2 = 17000      # K340 - Precio minimo
$ws.Cells.Item(340, 12).Value2 = 18000      # L340 - Precio maximo
$ws.Cells.Item(340, 13).Value2 = 17500      # M340 - Precio promedio ponderado
$ws.Cells.Item(340, 14).Value2 = "$/malla 10 kilos"  # N340 - Unidad de comercializacion
$ws.Cells.Item(340, 15).Value2 = "China"    # O340 - Origen
$ws.Cells.Item(340, 16).Value2 = 1750       # P340 - Precio $/Kg
$ws.Cells.Item(340, 17).Value2 = 10         # Q340 - Kg o Unidades

# 4) The brand-new row 442 needs the same date number format (style) that
#    column D carries throughout the rest of the table.
$ws.Range("D442").NumberFormat = $ws.Range("D441").NumberFormat
